{"js": "// Update the worksheet date and every math-problem cell in the single\n// 20x5 table. Problems are addressed by (row, col) position rather than\n// by matching old text, since several old values repeat (e.g. \"67-30=\",\n// \"20+49=\") at different positions with different replacements.\n\nconst NEW_DATE = \"2025-08-23 Saturday\";\n\n// GRID[row][col] holds the replacement text for the table cell at that\n// position (0-based), in the same order the table appears in the body.\nconst GRID = [\n  [\"87-84=\", \"35+10=\", \"86+9=\", \"1+30=\", \"49+6=\"],\n  [\"84-39=\", \"22-6=\", \"64+5=\", \"62-39=\", \"96-94=\"],\n  [\"82-13=\", \"85-73=\", \"15+7=\", \"72-18=\", \"29-19=\"],\n  [\"95-92=\", \"74-16=\", \"3+30=\", \"31+24=\", \"41+43=\"],\n  [\"81-37=\", \"27-18=\", \"39-24=\", \"19-17=\", \"69+7=\"],\n  [\"12+62=\", \"81+2=\", \"79+18=\", \"56-25=\", \"75-30=\"],\n  [\"19+71=\", \"93-86=\", \"93-32=\", \"63+4=\", \"28-1=\"],\n  [\"63-48=\", \"30+11=\", \"76-11=\", \"91-61=\", \"58-30=\"],\n  [\"24+57=\", \"41+58=\", \"52-45=\", \"87-40=\", \"85-38=\"],\n  [\"21+10=\", \"11+49=\", \"82-56=\", \"75-51=\", \"28+13=\"],\n  [\"42+1=\", \"73-72=\", \"50-34=\", \"7+68=\", \"13+11=\"],\n  [\"62-46=\", \"33+62=\", \"12+51=\", \"41-11=\", \"37+42=\"],\n  [\"18+72=\", \"22+39=\", \"46+6=\", \"20+79=\", \"67+23=\"],\n  [\"63-51=\", \"50+44=\", \"90-40=\", \"56-54=\", \"82+8=\"],\n  [\"67-2=\", \"19+45=\", \"22+20=\", \"25+25=\", \"28+47=\"],\n  [\"86-50=\", \"28+32=\", \"55+10=\", \"50-46=\", \"34-22=\"],\n  [\"19-18=\", \"75-13=\", \"15+36=\", \"42-17=\", \"56-26=\"],\n  [\"21-10=\", \"6+66=\", \"72-22=\", \"36+50=\", \"50+31=\"],\n  [\"79-55=\", \"52+36=\", \"56-29=\", \"21-13=\", \"46-7=\"],\n  [\"67-51=\", \"62-58=\", \"55-5=\", \"78-76=\", \"77+22=\"]\n];\n\n// 1) Update the date line (first paragraph in the document body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(NEW_DATE, Word.InsertLocation.replace);\n\n// 2) Update every cell of the (only) table, preserving per-cell formatting\n//    by setting TableCell.value instead of rewriting runs manually.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (let r = 0; r < GRID.length; r++) {\n  for (let c = 0; c < GRID[r].length; c++) {\n    table.getCell(r, c).value = GRID[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and every math-problem cell in the single\n# 20x5 table. Problems are addressed by (row, col) position rather than\n# by matching old text, since several old values repeat (e.g. \"67-30=\",\n# \"20+49=\") at different positions with different replacements.\n\n$newDate = '2025-08-23 Saturday'\n\n# $grid[row][col] holds the replacement text for the table cell at that\n# position (0-based), in the same order the table appears in the document.\n$grid = @(\n    @('87-84=', '35+10=', '86+9=', '1+30=', '49+6='),\n    @('84-39=', '22-6=', '64+5=', '62-39=', '96-94='),\n    @('82-13=', '85-73=', '15+7=', '72-18=', '29-19='),\n    @('95-92=', '74-16=', '3+30=', '31+24=', '41+43='),\n    @('81-37=', '27-18=', '39-24=', '19-17=', '69+7='),\n    @('12+62=', '81+2=', '79+18=', '56-25=', '75-30='),\n    @('19+71=', '93-86=', '93-32=', '63+4=', '28-1='),\n    @('63-48=', '30+11=', '76-11=', '91-61=', '58-30='),\n    @('24+57=', '41+58=', '52-45=', '87-40=', '85-38='),\n    @('21+10=', '11+49=', '82-56=', '75-51=', '28+13='),\n    @('42+1=', '73-72=', '50-34=', '7+68=', '13+11='),\n    @('62-46=', '33+62=', '12+51=', '41-11=', '37+42='),\n    @('18+72=', '22+39=', '46+6=', '20+79=', '67+23='),\n    @('63-51=', '50+44=', '90-40=', '56-54=', '82+8='),\n    @('67-2=', '19+45=', '22+20=', '25+25=', '28+47='),\n    @('86-50=', '28+32=', '55+10=', '50-46=', '34-22='),\n    @('19-18=', '75-13=', '15+36=', '42-17=', '56-26='),\n    @('21-10=', '6+66=', '72-22=', '36+50=', '50+31='),\n    @('79-55=', '52+36=', '56-29=', '21-13=', '46-7='),\n    @('67-51=', '62-58=', '55-5=', '78-76=', '77+22=')\n)\n\n$d = $word.ActiveDocument\n\n# 1) Update the date line (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = $newDate\n\n# 2) Update every cell of the (only) table, preserving per-cell formatting\n#    by writing straight into Cell.Range.Text (Word COM is 1-based).\n$t = $d.Tables.Item(1)\nfor ($r = 0; $r -lt $grid.Count; $r++) {\n    for ($c = 0; $c -lt $grid[$r].Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $grid[$r][$c]\n    }\n}\n"}
